$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBS")

# Add new row 4 with the "findEntryDateRangeFactorCustNo" lookup definition
$ws.Range("A4").Value = "findEntryDateRangeFactorCustNo"
$ws.Range("B4").Value = "EntryDate >= ,AND EntryDate <= ,AND Factor = , AND CustNo = "
$ws.Range("C4").Value = "EntryDate,DtlSeq ASC"

# Match the wrapped-text style used by the other rows in column B
$ws.Range("B4").WrapText = $true

# Move the active selection to C5, as recorded after the edit
$ws.Activate()
$ws.Range("C5").Select()
